$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3870.3867
$ws.Range("I64").Value = 3807.75
$ws.Range("J64").Value = 3928.205
$ws.Range("K64").Value = 3807.75
$ws.Range("L64").Value = 3928.205
$ws.Range("M64").Value = -3559.75
$ws.Range("N64").Value = -4424.205

$ws.Range("H67").Value = 3870.3867
$ws.Range("I67").Value = 3807.75
$ws.Range("J67").Value = 3928.205
$ws.Range("K67").Value = 3807.75
$ws.Range("L67").Value = 3928.205
$ws.Range("M67").Value = -2949.75
$ws.Range("N67").Value = -5644.205

$ws.Range("H74").Value = 3444.4443
$ws.Range("I74").Value = 3000
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -2064
$ws.Range("N74").Value = -5372

$ws.Range("H76").Value = 3444.6
$ws.Range("I76").Value = 2970.0425
$ws.Range("J76").Value = 4120.485
$ws.Range("K76").Value = 2970.0425
$ws.Range("L76").Value = 4120.485
$ws.Range("M76").Value = -2655.0425
$ws.Range("N76").Value = -4750.485

$ws.Range("H77").Value = 3444.4443
$ws.Range("I77").Value = 3000
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -10320
$ws.Range("N77").Value = -26860

$ws.Range("H79").Value = 3444.6
$ws.Range("I79").Value = 2970.0425
$ws.Range("J79").Value = 4120.485
$ws.Range("K79").Value = 2970.0425
$ws.Range("L79").Value = 4120.485
$ws.Range("M79").Value = -1878.0425
$ws.Range("N79").Value = -6304.485

$ws.Range("H82").Value = 6000106
$ws.Range("I82").Value = 6000106
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 18000318
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -17999912
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 6000106
$ws.Range("I85").Value = 6000106
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 18000318
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -17998914
$ws.Range("N85").ClearContents()

$ws.Range("H138").Value = 3122.8223
$ws.Range("I138").Value = 675.4737
$ws.Range("J138").Value = 4911.269
$ws.Range("K138").Value = 2026.4211
$ws.Range("L138").Value = 14733.807
$ws.Range("M138").Value = 3113.5789
$ws.Range("N138").Value = -25013.807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18522020
$ws.Range("I32").Value = 20001986
$ws.Range("J32").Value = 22456.75
$ws.Range("K32").Value = 20001986
$ws.Range("L32").Value = 22456.75
$ws.Range("M32").Value = -20001699
$ws.Range("N32").Value = -23030.75

$ws.Range("H61").Value = 3004.25
$ws.Range("I61").Value = 2561.4666
$ws.Range("J61").Value = 3742.2222
$ws.Range("K61").Value = 2561.4666
$ws.Range("L61").Value = 3742.2222
$ws.Range("M61").Value = -2349.4666
$ws.Range("N61").Value = -4166.2222

$ws.Range("H124").Value = 19975.572
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 19975.572
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 19975.572
$ws.Range("N124").Value = -29795.572

$ws.Range("H136").Value = 3004.25
$ws.Range("I136").Value = 2561.4666
$ws.Range("J136").Value = 3742.2222
$ws.Range("K136").Value = 7684.399800000001
$ws.Range("L136").Value = 11226.6666
$ws.Range("M136").Value = -5134.399800000001
$ws.Range("N136").Value = -16326.6666

$ws.Range("H140").Value = 44000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 44000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 44000
$ws.Range("N140").Value = -54360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H82").Value = 5689
$ws.Range("I82").Value = 3427.9
$ws.Range("J82").Value = 28300
$ws.Range("K82").Value = 3427.9
$ws.Range("L82").Value = 28300
$ws.Range("M82").Value = -3044.9
$ws.Range("N82").Value = -29066

$ws.Range("H85").Value = 5689
$ws.Range("I85").Value = 3427.9
$ws.Range("J85").Value = 28300
$ws.Range("K85").Value = 3427.9
$ws.Range("L85").Value = 28300
$ws.Range("M85").Value = -2101.9
$ws.Range("N85").Value = -30952

$ws.Range("H134").Value = 1850.3715
$ws.Range("I134").Value = 1795.8485
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 5387.5455
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -2852.5455
$ws.Range("N134").Value = -13320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1579.2941
$ws.Range("I132").Value = 963
$ws.Range("J132").Value = 3956.4285
$ws.Range("K132").Value = 2889
$ws.Range("L132").Value = 11869.2855
$ws.Range("M132").Value = -359
$ws.Range("N132").Value = -16929.2855

$ws.Range("H134").Value = 4685.4443
$ws.Range("I134").Value = 1278.5
$ws.Range("J134").Value = 7411
$ws.Range("K134").Value = 3835.5
$ws.Range("L134").Value = 22233
$ws.Range("M134").Value = -1300.5
$ws.Range("N134").Value = -27303

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 555795.4
$ws.Range("I92").Value = 1000162.3
$ws.Range("J92").Value = 336.75
$ws.Range("K92").Value = 3000486.9
$ws.Range("L92").Value = 1010.25
$ws.Range("M92").Value = -2999238.9
$ws.Range("N92").Value = -3506.25

$ws.Range("H121").Value = 970.34485
$ws.Range("I121").Value = 483.33334
$ws.Range("J121").Value = 1026.5385
$ws.Range("K121").Value = 1450.00002
$ws.Range("L121").Value = 3079.6155
$ws.Range("M121").Value = -140.0000199999999
$ws.Range("N121").Value = -5699.6155

$ws.Range("H131").Value = 794.75
$ws.Range("I131").Value = 370
$ws.Range("J131").Value = 955.86206
$ws.Range("K131").Value = 1110
$ws.Range("L131").Value = 2867.58618
$ws.Range("M131").Value = 3930
$ws.Range("N131").Value = -12947.58618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5167.551
$ws.Range("I70").Value = 5348.769
$ws.Range("J70").Value = 4460.8
$ws.Range("K70").Value = 5348.769
$ws.Range("L70").Value = 4460.8
$ws.Range("M70").Value = -5078.769
$ws.Range("N70").Value = -5000.8

$ws.Range("H73").Value = 5167.551
$ws.Range("I73").Value = 5348.769
$ws.Range("J73").Value = 4460.8
$ws.Range("K73").Value = 5348.769
$ws.Range("L73").Value = 4460.8
$ws.Range("M73").Value = -4412.769
$ws.Range("N73").Value = -6332.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2032.4166
$ws.Range("I7").Value = 2086.75
$ws.Range("J7").Value = 1923.75
$ws.Range("K7").Value = 2086.75
$ws.Range("L7").Value = 1923.75
$ws.Range("M7").Value = -1974.75
$ws.Range("N7").Value = -2147.75

$ws.Range("H125").Value = 49905
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 49905
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 49905
$ws.Range("N125").Value = -59745

$ws.Range("H126").Value = 2032.4166
$ws.Range("I126").Value = 2086.75
$ws.Range("J126").Value = 1923.75
$ws.Range("K126").Value = 6260.25
$ws.Range("L126").Value = 5771.25
$ws.Range("M126").Value = -3790.25
$ws.Range("N126").Value = -10711.25

$ws.Range("H127").Value = 39333.332
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 39333.332
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 39333.332
$ws.Range("N127").Value = -49253.332

$ws.Range("H132").Value = 2615.4482
$ws.Range("I132").Value = 2457
$ws.Range("J132").Value = 2874.7273
$ws.Range("K132").Value = 7371
$ws.Range("L132").Value = 8624.1819
$ws.Range("M132").Value = -4841
$ws.Range("N132").Value = -13684.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1053
$ws.Range("I81").Value = 1082.8572
$ws.Range("J81").Value = 983.3333
$ws.Range("K81").Value = 2165.7144
$ws.Range("L81").Value = 1966.6666
$ws.Range("M81").Value = -1104.7144
$ws.Range("N81").Value = -4088.6666

$ws.Range("H84").Value = 1053
$ws.Range("I84").Value = 1082.8572
$ws.Range("J84").Value = 983.3333
$ws.Range("K84").Value = 10828.572
$ws.Range("L84").Value = 9833.333
$ws.Range("M84").Value = -5524.572
$ws.Range("N84").Value = -20441.333
